$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while guaranteeing it is stored as literal text,
# since Excel would otherwise silently reinterpret numeric-looking strings
# (e.g. "10.00", "1.00", "0.176") as numbers and lose their exact formatting.
# We temporarily force Text number format for the assignment, then restore
# the cell's original style so no visual/style differences are introduced.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "72.746.89"
Set-TextValue "E2" "  +4.10%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.636.36"
Set-TextValue "E3" "  +2.66%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "605.63"
Set-TextValue "E5" "  +0.90%  "

# Row 6 - Solana
Set-TextValue "D6" "179.22"
Set-TextValue "E6" "  +0.41%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.176"
Set-TextValue "E9" "  +9.25%  "

# Row 10 - LidoStakedEther
Set-TextValue "D10" "2.635.61"
Set-TextValue "E10" "  +2.62%  "

# Row 11 - TRON
Set-TextValue "E11" "  +1.29%  "

# Row 12 - Cardano
Set-TextValue "E12" "  +3.18%  "

# Row 13 - Toncoin
Set-TextValue "E13" "  +0.18%  "

# Row 14 - ShibaInu
Set-TextValue "D14" "0.0000191"
Set-TextValue "E14" "  +4.58%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "E15" "  +3.21%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "72.556.63"
Set-TextValue "E16" "  +4.04%  "

# Row 17 - Avalanche
Set-TextValue "E17" "  +2.00%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.627.49"
Set-TextValue "E18" "  +3.21%  "

# Row 19 - now BitcoinCash (was Chainlink)
Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "385.64"
Set-TextValue "E19" "  +5.41%  "

# Row 20 - now Chainlink (was BitcoinCash)
Set-TextValue "B20" "Chainlink"
Set-TextValue "C20" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "11.76"
Set-TextValue "E20" "  +5.03%  "

# Row 21 - Uniswap
Set-TextValue "E21" "  +1.88%  "

# Row 22 - Polkadot
Set-TextValue "E22" "  +1.30%  "

# Row 23 - SuiNetwork
Set-TextValue "E23" "  +15.86%  "

# Row 24 - Litecoin
Set-TextValue "D24" "74.19"

# Row 26 - Dai
Set-TextValue "E26" "  +0.09%  "

# Row 27 - Aptos
Set-TextValue "D27" "10.00"
Set-TextValue "E27" "  +8.41%  "

# Row 28 - WrappedeETH
Set-TextValue "D28" "2.730.22"

# Row 29 - Binance-PegBSC-USD
Set-TextValue "E29" "  +0.05%  "

# Row 30 - PEPE
Set-TextValue "D30" "0.0₃0962"
Set-TextValue "E30" "  +4.66%  "

# Row 31 - now InternetComputer(DFINITY) (was Bittensor)
Set-TextValue "B31" "InternetComputer(DFINITY)"
Set-TextValue "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "8.12"
Set-TextValue "E31" "  +3.94%  "

# Row 32 - now Bittensor (was InternetComputer(DFINITY))
Set-TextValue "B32" "Bittensor"
Set-TextValue "C32" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D32" "519.26"
Set-TextValue "E32" "  +0.91%  "

# Row 33 - Fetch.AI
Set-TextValue "E33" "  +4.30%  "

# Row 34 - PancakeSwap
Set-TextValue "E34" "  +1.61%  "

# Row 35 - FirstDigitalUSD
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  -0.07%  "

# Row 36 - Monero
Set-TextValue "D36" "162.99"
Set-TextValue "E36" "  -0.65%  "

# Row 37 - EthereumClassic
Set-TextValue "D37" "19.42"
Set-TextValue "E37" "  +2.25%  "

# Row 38 - ImmutableX
Set-TextValue "E38" "  +3.95%  "

# Row 39 - now Kaspa (was WhiteBITCoin)
Set-TextValue "B39" "Kaspa"
Set-TextValue "C39" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.112"
Set-TextValue "E39" "  -5.77%  "

# Row 40 - now WhiteBITCoin (was Kaspa)
Set-TextValue "B40" "WhiteBITCoin"
Set-TextValue "C40" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D40" "19.12"
Set-TextValue "E40" "  +1.09%  "

# Row 41 - Stacks
Set-TextValue "E41" "  +5.54%  "

# Row 42 - RenderToken
Set-TextValue "E42" "  +4.61%  "

# Row 43 - USDe
Set-TextValue "E43" "  +0.07%  "

# Row 44 - dogwifhat
Set-TextValue "D44" "2.60"
Set-TextValue "E44" "  +5.14%  "

# Row 45 - PolygonEcosystemToken
Set-TextValue "E45" "  +2.74%  "

# Row 46 - OKB
Set-TextValue "E46" "  +1.05%  "

# Row 47 - Aave
Set-TextValue "D47" "151.22"
Set-TextValue "E47" "  -0.45%  "

# Row 48 - Filecoin
Set-TextValue "E48" "  +1.94%  "

# Row 49 - ARBITRUM
Set-TextValue "E49" "  +4.44%  "

# Row 50 - Optimism
Set-TextValue "E50" "  +4.95%  "

# Row 51 - BabyDogeCoin
Set-TextValue "D51" "0.0₆0265"
Set-TextValue "E51" "  +2.19%  "
